$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44330, 0, 4, 87.24100327153762),
    @(44331, 1, 5, 109.051254089422),
    @(44332, 0, 4, 87.24100327153762),
    @(44333, 0, 3, 65.43075245365321),
    @(44334, 0, 3, 65.43075245365321),
    @(44335, 0, 3, 65.43075245365321),
    @(44336, 0, 1, 21.81025081788441),
    @(44337, 0, 1, 21.81025081788441),
    @(44338, 1, 1, 21.81025081788441),
    @(44339, 2, 3, 65.43075245365321),
    @(44340, 0, 3, 65.43075245365321),
    @(44341, 0, 3, 65.43075245365321),
    @(44342, 0, 3, 65.43075245365321),
    @(44343, 1, 4, 87.24100327153762)
)

$lastExistingRow = 255
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the date-column formatting/style from the previous row so the
    # new dates keep the same style (same as column A of row 255).
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
